# Updated cryptos list on Sat Apr  6 12:39:02 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    # The Price column stores plain text (e.g. "580.71"), but Excel's COM
    # layer auto-coerces strings that parse cleanly as numbers into real
    # numeric cells. Force a Text number format while writing the value,
    # then restore the default "Normal" style so the cell ends up looking
    # exactly like the untouched ones (no explicit style / number format).
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue "D2" "67.702.77"
$ws.Range("E2").Value = "  +2.35%  "

# Row 3 - Ethereum
Set-TextValue "D3" "3.332.81"
$ws.Range("E3").Value = "  +3.32%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.15%  "

# Row 5 - BNB
Set-TextValue "D5" "580.48"
$ws.Range("E5").Value = "  +0.87%  "

# Row 6 - Solana
Set-TextValue "D6" "175.59"
$ws.Range("E6").Value = "  +4.13%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  -0.14%  "

# Row 8 - XRP
$ws.Range("E8").Value = "  +2.92%  "

# Row 9 - LidoStakedEther
Set-TextValue "D9" "3.328.56"
$ws.Range("E9").Value = "  +3.40%  "

# Row 10 - Dogecoin
Set-TextValue "D10" "0.182"
$ws.Range("E10").Value = "  +8.08%  "

# Row 11 - Cardano
$ws.Range("E11").Value = "  +3.16%  "

# Row 12 - Avalanche
Set-TextValue "D12" "46.80"
$ws.Range("E12").Value = "  +6.11%  "

# Row 13 - ShibaInu
$ws.Range("E13").Value = "  +2.14%  "

# Row 14 - BitcoinCash
Set-TextValue "D14" "690.70"
$ws.Range("E14").Value = "  +2.94%  "

# Row 15 - WrappedliquidstakedEther2.0
Set-TextValue "D15" "3.873.87"
$ws.Range("E15").Value = "  +3.14%  "

# Row 16 - Polkadot
$ws.Range("E16").Value = "  +3.24%  "

# Row 17 - WrappedBTC
Set-TextValue "D17" "67.688.95"
$ws.Range("E17").Value = "  +2.01%  "

# Row 18 - TRON
$ws.Range("E18").Value = "  +0.94%  "

# Row 19 - WrappedEther
Set-TextValue "D19" "3.325.80"
$ws.Range("E19").Value = "  +2.75%  "

# Row 20 - Chainlink
Set-TextValue "D20" "17.55"
$ws.Range("E20").Value = "  +3.40%  "

# Row 21 - Uniswap
Set-TextValue "D21" "11.05"
$ws.Range("E21").Value = "  +4.77%  "

# Row 22 - Polygon
Set-TextValue "D22" "0.892"
$ws.Range("E22").Value = "  +3.27%  "

# Row 23 - Toncoin
Set-TextValue "D23" "5.53"
$ws.Range("E23").Value = "  +6.67%  "

# Row 24 - InternetComputer(DFINITY)
Set-TextValue "D24" "16.85"
$ws.Range("E24").Value = "  +0.80%  "

# Row 25 - Litecoin
Set-TextValue "D25" "100.94"
$ws.Range("E25").Value = "  +5.59%  "

# Row 26 - PancakeSwap
$ws.Range("E26").Value = "  +2.78%  "

# Row 27 - ImmutableX
$ws.Range("E27").Value = "  +3.06%  "

# Row 28 - RenderToken
Set-TextValue "D28" "9.39"
$ws.Range("E28").Value = "  +5.95%  "

# Row 29 - EthereumClassic
Set-TextValue "D29" "32.98"
$ws.Range("E29").Value = "  +3.58%  "

# Row 30 - Filecoin
Set-TextValue "D30" "8.53"
$ws.Range("E30").Value = "  +4.79%  "

# Row 31 - NEARProtocol
Set-TextValue "D31" "7.03"
$ws.Range("E31").Value = "  +6.52%  "

# Row 32 - Bittensor
Set-TextValue "D32" "567.60"
$ws.Range("E32").Value = "  +0.60%  "

# Row 33 - Cosmos
Set-TextValue "D33" "11.01"
$ws.Range("E33").Value = "  +2.79%  "

# Row 34 - Hedera
$ws.Range("E34").Value = "  +4.44%  "

# Row 35 - OKB
Set-TextValue "D35" "57.45"
$ws.Range("E35").Value = "  +4.95%  "

# Row 36 - was Maker, now Dai
$ws.Range("B36").Value = "Dai"
$ws.Range("C36").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextValue "D36" "1.00"
$ws.Range("E36").Value = "  -0.02%  "

# Row 37 - was Dai, now Maker
$ws.Range("B37").Value = "Maker"
$ws.Range("C37").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextValue "D37" "3.707.69"
$ws.Range("E37").Value = "  -2.03%  "

# Row 38 - dogwifhat
Set-TextValue "D38" "3.26"
$ws.Range("E38").Value = "  -0.55%  "

# Row 39 - InjectiveProtocol
Set-TextValue "D39" "35.09"
$ws.Range("E39").Value = "  +12.75%  "

# Row 40 - Kaspa
$ws.Range("E40").Value = "  +4.53%  "

# Row 41 - Stacks
Set-TextValue "D41" "3.15"
$ws.Range("E41").Value = "  +7.07%  "

# Row 42 - Fetch.AI
$ws.Range("E42").Value = "  +3.14%  "

# Row 43 - ApeXProtocol
Set-TextValue "D43" "3.34"
$ws.Range("E43").Value = "  +1.03%  "

# Row 44 - TheGraph
Set-TextValue "D44" "0.334"
$ws.Range("E44").Value = "  +5.05%  "

# Row 45 - PEPE
$ws.Range("E45").Value = "  +4.71%  "

# Row 46 - VeChain
Set-TextValue "D46" "0.0406"
$ws.Range("E46").Value = "  +3.73%  "

# Row 47 - ThetaToken
Set-TextValue "D47" "2.64"
$ws.Range("E47").Value = "  +6.76%  "

# Row 48 - Stellar
$ws.Range("E48").Value = "  +3.02%  "

# Row 49 - FirstDigitalUSD
$ws.Range("E49").Value = "  -0.33%  "

# Row 50 - Mantle
$ws.Range("E50").Value = "  +2.10%  "

# Row 51 - Monero
Set-TextValue "D51" "132.23"
$ws.Range("E51").Value = "  +4.24%  "
